# Add new "Release 3" task rows (25-29) to the HR Portal tracking sheet,
# along with the new comments/notes these tasks introduce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 (SL No 25) - height 63 like other long-comment rows (e.g. row 11)
$ws.Cells.Item(27, 1).Value2 = 25
$ws.Cells.Item(27, 2).Value2 = "Bonus Column should add in payroll and should show in payslip if bonus present"
$ws.Cells.Item(27, 3).Value2 = "Backend & UI"
$ws.Cells.Item(27, 4).Value2 = "Completed"
$ws.Cells.Item(27, 5).Value2 = "Changes in PayrollDO, jsp and Util"
$ws.Rows.Item(27).RowHeight = 63

# Row 28 (SL No 26)
$ws.Cells.Item(28, 1).Value2 = 26
$ws.Cells.Item(28, 2).Value2 = "Edit Access to HRUser also"
$ws.Cells.Item(28, 3).Value2 = "UI"
$ws.Cells.Item(28, 4).Value2 = "Completed"
$ws.Cells.Item(28, 5).Value2 = "Changes in viewemployee.jsp"
$ws.Rows.Item(28).RowHeight = 31.5

# Row 29 (SL No 27)
$ws.Cells.Item(29, 1).Value2 = 27
$ws.Cells.Item(29, 2).Value2 = "Experience input field change to year and month"
$ws.Cells.Item(29, 3).Value2 = "UI"
$ws.Cells.Item(29, 4).Value2 = "Completed"
$ws.Cells.Item(29, 5).Value2 = "Changes in viewemployee.jsp & addemployee.jsp"
$ws.Rows.Item(29).RowHeight = 31.5

# Row 30 (SL No 28)
$ws.Cells.Item(30, 1).Value2 = 28
$ws.Cells.Item(30, 2).Value2 = "Update not working for middle name and languages known in employee"
$ws.Cells.Item(30, 3).Value2 = "UI"
$ws.Cells.Item(30, 4).Value2 = "Completed"
$ws.Cells.Item(30, 5).Value2 = "Changes in viewemployee.jsp"
$ws.Rows.Item(30).RowHeight = 63

# Row 31 (SL No 29)
$ws.Cells.Item(31, 1).Value2 = 29
$ws.Cells.Item(31, 2).Value2 = "view employee  arrow images and showing same tab when again come to that page"
$ws.Cells.Item(31, 3).Value2 = "UI"
$ws.Cells.Item(31, 4).Value2 = "Completed"
$ws.Cells.Item(31, 5).Value2 = "Changes in viewemployee.jsp"
$ws.Rows.Item(31).RowHeight = 63

# Update the view/selection to mirror the author's scroll position & active cell
$ws.Range("E31").Select()
